$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '257.91'
$r.Style = "Normal"
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '22.77'
$r.Style = "Normal"
$r = $ws.Range('D4')
$r.NumberFormat = "@"
$r.Value = '6.150'
$r.Style = "Normal"
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '0.06078'
$r.Style = "Normal"
$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '3.446'
$r.Style = "Normal"
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '1.363'
$r.Style = "Normal"
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.7993'
$r.Style = "Normal"
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '0.01321'
$r.Style = "Normal"
$ws.Range('E10').Value = '9OneONE'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '0.1581'
$r.Style = "Normal"
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '0.08036'
$r.Style = "Normal"
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '0.03344'
$r.Style = "Normal"
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '0.03100'
$r.Style = "Normal"
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '0.09299'
$r.Style = "Normal"
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '3.898'
$r.Style = "Normal"
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '0.001700'
$r.Style = "Normal"
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '0.04839'
$r.Style = "Normal"
$ws.Range('E18').Value = '17CoinExTokenCET'
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '0.006252'
$r.Style = "Normal"
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '0.001102'
$r.Style = "Normal"
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '0.003385'
$r.Style = "Normal"
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '0.0001500'
$r.Style = "Normal"
$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '3.686'
$r.Style = "Normal"
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '2.264'
$r.Style = "Normal"
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '0.3358'
$r.Style = "Normal"
$r = $ws.Range('D26')
$r.NumberFormat = "@"
$r.Value = '0.1272'
$r.Style = "Normal"
$r = $ws.Range('D27')
$r.NumberFormat = "@"
$r.Value = '0.0003016'
$r.Style = "Normal"
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '0.007162'
$r.Style = "Normal"
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '0.003901'
$r.Style = "Normal"
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '0.1112'
$r.Style = "Normal"
$ws.Range('E43').Value = '42BKEXTokenBKK'
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '0.009977'
$r.Style = "Normal"
$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '0.00005927'
$r.Style = "Normal"
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '0.00000000750'
$r.Style = "Normal"
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '0.7501'
$r.Style = "Normal"
$r = $ws.Range('D49')
$r.NumberFormat = "@"
$r.Value = '0.07120'
$r.Style = "Normal"
$ws.Range('E49').Value = '48BOLOBOLOWorstin24h'
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '0.00001500'
$r.Style = "Normal"
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '0.01010'
$r.Style = "Normal"
